# Add a styled header row (monster-stat-block column headers) to Лист1 (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column A : "Key" -> Check Cell style (gray fill, double border, bold white custom font) ---
$ws.Range("A1").Value = "Key"
$ws.Range("A1").Style = "Check Cell"
$ws.Range("A1").Font.Name = "18thCentury"
$ws.Range("A1").Borders.LineStyle = -4119   # xlDouble

# --- Column B : "Name" -> Heading 1 style, red ---
$ws.Range("B1").Value = "Name"
$ws.Range("B1").Style = "Heading 1"
$ws.Range("B1").Font.Color = 255            # FF0000 red

# --- Columns C:I : Type, Size, Alignment, Environment, AC, HP, Speed -> Heading 1 style, green ---
$greenCols = @("C","D","E","F","G","H","I")
$greenVals = @("Type","Size","Alignment","Environment","AC","HP","Speed")
for ($i = 0; $i -lt $greenCols.Length; $i++) {
    $cell = $ws.Range($greenCols[$i] + "1")
    $cell.Value = $greenVals[$i]
    $cell.Style = "Heading 1"
    $cell.Font.Color = 5287936              # 00B050 green
}

# --- Columns J:O : STR, DEX, CON, INT, WIS, CHA -> Heading 1 style, blue ---
$blueCols = @("J","K","L","M","N","O")
$blueVals = @("STR","DEX","CON","INT","WIS","CHA")
for ($i = 0; $i -lt $blueCols.Length; $i++) {
    $cell = $ws.Range($blueCols[$i] + "1")
    $cell.Value = $blueVals[$i]
    $cell.Style = "Heading 1"
    $cell.Font.Color = 15773696             # 00B0F0 blue
}

# --- Columns P:V : Saving Throws, Skills, Damage immunities, Condition immunities, Senses, Languages, Challenge -> Heading 1 style, orange ---
$orangeCols = @("P","Q","R","S","T","U","V")
$orangeVals = @("Saving Throws","Skills","Damage immunities","Condition immunities","Senses","Languages","Challenge")
for ($i = 0; $i -lt $orangeCols.Length; $i++) {
    $cell = $ws.Range($orangeCols[$i] + "1")
    $cell.Value = $orangeVals[$i]
    $cell.Style = "Heading 1"
    $cell.Font.Color = 49407                # FFC000 orange
}

# --- Column W : "Description" -> Heading 1 style, purple ---
$ws.Range("W1").Value = "Description"
$ws.Range("W1").Style = "Heading 1"
$ws.Range("W1").Font.Color = 10498160       # 7030A0 purple

# Row heights (header row taller; row below kept for the thick bottom border look)
$ws.Rows.Item(1).RowHeight = 21
$ws.Rows.Item(2).RowHeight = 15.75

# Page setup: portrait, paper size 9 (A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final selection on the last header cell, like the source file
[void]$ws.Range("W1").Select()

Write-Host "Header row written to Лист1"
